$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    # Force Excel to store the value as literal text (it would otherwise
    # auto-convert plain-numeric-looking strings to numbers), then restore
    # the cell to the default "Normal" style so no formatting is changed.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "28.014.58"
$ws.Range("E2").Value = "  +7.00%  "

$ws.Range("D3").Value = "1.742.21"
$ws.Range("E3").Value = "  +5.30%  "

Set-TextValue $ws.Range("D4") "1.003"
$ws.Range("E4").Value = "  -0.19%  "

Set-TextValue $ws.Range("D5") "228.49"
$ws.Range("E5").Value = "  +4.16%  "

Set-TextValue $ws.Range("D6") "0.5432"
$ws.Range("E6").Value = "  +3.39%  "

$ws.Range("E7").Value = "  -0.25%  "

Set-TextValue $ws.Range("D8") "0.2766"
$ws.Range("E8").Value = "  +3.67%  "

Set-TextValue $ws.Range("D9") "0.06723"
$ws.Range("E9").Value = "  +5.78%  "

Set-TextValue $ws.Range("D10") "21.72"
$ws.Range("E10").Value = "  +4.92%  "

Set-TextValue $ws.Range("D11") "0.07795"
$ws.Range("E11").Value = "  +0.91%  "

Set-TextValue $ws.Range("D12") "4.699"
$ws.Range("E12").Value = "  +2.08%  "

$ws.Range("D13").Value = "1.749.88"
$ws.Range("E13").Value = "  +2.84%  "

$ws.Range("D14").Value = "1.982.20"
$ws.Range("E14").Value = "  +5.23%  "

Set-TextValue $ws.Range("D15") "0.5967"
$ws.Range("E15").Value = "  +6.00%  "

$ws.Range("D16").Value = "0.0₅8381"
$ws.Range("E16").Value = "  +2.10%  "

Set-TextValue $ws.Range("D17") "68.87"
$ws.Range("E17").Value = "  +5.22%  "

$ws.Range("D18").Value = "28.009.78"
$ws.Range("E18").Value = "  +6.95%  "

Set-TextValue $ws.Range("D19") "225.71"
$ws.Range("E19").Value = "  +17.52%  "

Set-TextValue $ws.Range("D20") "4.839"
$ws.Range("E20").Value = "  +2.69%  "

$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("E22").Value = "  +4.87%  "

Set-TextValue $ws.Range("D23") "6.231"
$ws.Range("E23").Value = "  +3.91%  "

Set-TextValue $ws.Range("D24") "1.003"
$ws.Range("E24").Value = "  -0.22%  "

Set-TextValue $ws.Range("D25") "146.28"
$ws.Range("E25").Value = "  +1.25%  "

Set-TextValue $ws.Range("D26") "0.1249"
$ws.Range("E26").Value = "  +3.37%  "

Set-TextValue $ws.Range("D27") "7.473"
$ws.Range("E27").Value = "  +2.71%  "

Set-TextValue $ws.Range("D28") "17.17"
$ws.Range("E28").Value = "  +7.54%  "

Set-TextValue $ws.Range("D29") "1.644"
$ws.Range("E29").Value = "  +9.88%  "

Set-TextValue $ws.Range("D30") "0.05663"
$ws.Range("E30").Value = "  +0.54%  "

Set-TextValue $ws.Range("D31") "1.317"
$ws.Range("E31").Value = "  +3.35%  "

Set-TextValue $ws.Range("D32") "3.704"
$ws.Range("E32").Value = "  +5.65%  "

Set-TextValue $ws.Range("D33") "3.517"
$ws.Range("E33").Value = "  +4.09%  "

Set-TextValue $ws.Range("D34") "1.676"
$ws.Range("E34").Value = "  +5.86%  "

Set-TextValue $ws.Range("D35") "0.9826"
$ws.Range("E35").Value = "  +3.08%  "

$ws.Range("E36").Value = "  +2.24%  "

Set-TextValue $ws.Range("D37") "2.451"
$ws.Range("E37").Value = "  +1.70%  "

Set-TextValue $ws.Range("D38") "0.5956"
$ws.Range("E38").Value = "  +3.30%  "

Set-TextValue $ws.Range("D39") "0.01667"
$ws.Range("E39").Value = "  +4.32%  "

Set-TextValue $ws.Range("D40") "5.953"
$ws.Range("E40").Value = "  -0.93%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D41") "0.8496"
$ws.Range("E41").Value = "  +0.90%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.050.35"
$ws.Range("E42").Value = "  +3.62%  "

Set-TextValue $ws.Range("D44") "102.07"
$ws.Range("E44").Value = "  +0.17%  "

$ws.Range("D45").Value = "1.886.90"
$ws.Range("E45").Value = "  +5.08%  "

$ws.Range("E46").Value = "  +4.78%  "

Set-TextValue $ws.Range("D47") "59.99"
$ws.Range("E47").Value = "  +2.52%  "

Set-TextValue $ws.Range("D48") "8.270"
$ws.Range("E48").Value = "  +2.98%  "

Set-TextValue $ws.Range("D49") "0.4433"
$ws.Range("E49").Value = "  +1.84%  "

Set-TextValue $ws.Range("D50") "0.05325"
$ws.Range("E50").Value = "  -0.37%  "

Set-TextValue $ws.Range("D51") "0.9960"
$ws.Range("E51").Value = "  -0.88%  "

